# "updated activity till excel form"
# Row 2 (first innings entry) loses its runs/balls/fours, row 4 (the
# previously-blank entry) picks them up instead: 6/7/1 moves from row 2 to row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these as text cells (matching the workbook's numberStoredAsText
# convention) instead of letting Excel auto-coerce numeric-looking strings
# into real numbers.
$ws.Range("C2:E2").NumberFormat = "@"
$ws.Range("C4:E4").NumberFormat = "@"

$ws.Range("C2").Value = "0"
$ws.Range("D2").Value = "0"
$ws.Range("E2").Value = "0"

$ws.Range("C4").Value = "6"
$ws.Range("D4").Value = "7"
$ws.Range("E4").Value = "1"
